# "made great progress with add product page"
# A new functionality row is appended to the tracking sheet: "Update Product",
# still unsolved (red cell, same style as the other not-yet-done rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 25: "Update Product" functionality, not yet solved ---
# Duplicate the formatting of an existing "unsolved" row (row 23: red/unsolved
# style in column C, thick-bordered row) onto the new row 25, then overwrite
# the cell values for the new entry.
$ws.Range("A23:C23").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Update Product"
$ws.Range("C25").ClearContents()

# --- Update the active selection to reflect where the user left off ---
$ws.Range("D9").Select()
